# Generate Report for Archive
#
# 1) The localization status moved on from "Ready for handoff" to
#    "In Translation" for both tracked files, on every sheet that surfaces
#    the status column (Overview!E:F, zh-cn!C, de-de!C).
# 2) Those status columns got narrower once the shorter text was in place
#    (columns were kept auto-sized to their contents).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status shown in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 5).Value = "In Translation"
$wsOverview.Cells.Item(2, 6).Value = "In Translation"
$wsOverview.Cells.Item(3, 5).Value = "In Translation"
$wsOverview.Cells.Item(3, 6).Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn detail sheet: status shown in column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(2, 3).Value = "In Translation"
$wsZhCn.Cells.Item(3, 3).Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: status shown in column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(2, 3).Value = "In Translation"
$wsDeDe.Cells.Item(3, 3).Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
